$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row9 = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
$row10 = @(10, 9, 8, 7, 6, 5, 4, 3, 2, 1)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $row9[$i]
    $ws.Cells.Item(10, $i + 1).Value = $row10[$i]
}

$ws.Range("J10").Select()
